# Take the calculation of some variables out of the loop
$wb = $excel.ActiveWorkbook

$wsRunControl = $wb.Worksheets.Item("RunControl")
$wsGlobalParams = $wb.Worksheets.Item("GlobalParams")

# Set K14:K19 and K22:K27 ("no new entrants" switch) to FALSE since the
# variable calculation was pulled out of the loop (was constant TRUE).
$wsRunControl.Range("K14:K19").Value = $false
$wsRunControl.Range("K22:K27").Value = $false

# Update GlobalParams: nyear 105 -> 80, nsim 1000 -> 10000
$wsGlobalParams.Range("A3").Value = 80
$wsGlobalParams.Range("B3").Value = 10000

# Update selections / active cells to match the new active sheet state
$wsRunControl.Range("K29").Select()
$wsGlobalParams.Range("B3").Select()

# Make GlobalParams the active sheet (was RunControl)
$wsGlobalParams.Activate()
